$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 282, shifting existing rows 282:344 down to 283:345.
$ws.Rows.Item(282).Insert()

# Populate the newly inserted row 282 with the new record.
$ws.Cells.Item(282, 1).Value2 = 5
$ws.Cells.Item(282, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(282, 3).Value2 = "Maule"
$ws.Cells.Item(282, 4).Value2 = 44711
$ws.Cells.Item(282, 5).Value2 = 7
$ws.Cells.Item(282, 6).Value2 = 100114013
$ws.Cells.Item(282, 7).Value2 = "Zanahoria"
$ws.Cells.Item(282, 8).Value2 = "Sin especificar"
$ws.Cells.Item(282, 9).Value2 = "Primera"
$ws.Cells.Item(282, 10).Value2 = 500
$ws.Cells.Item(282, 11).Value2 = 5500
$ws.Cells.Item(282, 12).Value2 = 5500
$ws.Cells.Item(282, 13).Value2 = 5500
$ws.Cells.Item(282, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(282, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(282, 16).Value2 = 275
$ws.Cells.Item(282, 17).Value2 = 20
$ws.Cells.Item(282, 18).Value2 = "Hortaliza"

# Apply the same date number format used by the other rows in column D.
$ws.Cells.Item(282, 4).NumberFormat = $ws.Cells.Item(283, 4).NumberFormat
